# AFDP-3458: Add new MyDocuments module
# - Deny read access to * participant for PERSONAL document repositories.
#
# This inserts a new "Document Repository -deny read access" rule row right
# after the existing "Document Repository -default read access" row (the
# DOC_REPO rule block starts at row 59), pushing the remaining DOC_REPO rule
# rows (Anybody can add comments / Lockout No Access Users / Anybody can add
# tag / Anybody can subscribe / Restricted Flag / Only participants can add
# files / Only participants can save / Only participants can upload or
# replace files) down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new blank row at 60; rows 60:67 shift down to 61:68.
$ws.Rows("60:60").Insert()

# Row-insert copies formatting from the row above into every column of the
# new row, including column A (which is otherwise unused below row 59) -
# drop that stray cell so the new row matches the other DOC_REPO rule rows.
$ws.Range("A60").Clear()

# Restore the correct per-column look-and-feel for the new row: columns B
# and G use the "header-ish" rule/outcome style, columns C:F use the plain
# condition-cell style - copy those styles (not values) from existing cells
# that already use them.
$ws.Range("B61").Copy()
$ws.Range("B60").PasteSpecial(-4122)
$ws.Range("G61").Copy()
$ws.Range("G60").PasteSpecial(-4122)
$ws.Range("C59:F59").Copy()
$ws.Range("C60:F60").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Match the row height used by the other short (single-expression) rule rows.
$ws.Rows("60:60").RowHeight = 30

# New rule: deny read access to "*" for repositories whose repositoryType is
# PERSONAL (mandatory deny read, so it cannot be overridden by other rules).
$ws.Range("B60").Value = "Document Repository -deny read access"
$ws.Range("C60").Value = "DOC_REPO"
$ws.Range("D60").Value = "repositoryType == 'PERSONAL'"
$ws.Range("G60").Value = "mandatory deny read to *"

# Leave the selection where the author's saved view landed.
[void]$ws.Range("I59").Select()
